$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "TipoDeHabilidade" side-tables (originally at G8:I12 and K8:L12, to the
# right of the "Habilidade"/"ClasseHabilidade" tables) are relocated below
# the existing tables, starting two rows under the last one (row 16 -> gap
# at row 18 -> new tables start at row 19).
#
# Column layout changes too: the "idTipo/NomeTipo" block (old K:L) becomes
# the new A:B block, and the "DeuBug" block (old G:I) becomes the new D:F
# block.

# Move K8:L12 ("TipoDeHabilidade" / idTipo / NomeTipo table) to A19:B23
$ws.Range("K8:L12").Cut($ws.Range("A19"))

# Move G8:I12 ("DeuBug" table) to D19:F23
$ws.Range("G8:I12").Cut($ws.Range("D19"))

# The old locations (columns G:L, rows 8-12) are now completely empty;
# remove the leftover (now valueless) cell formatting so no stray <c>
# elements remain there.
$ws.Range("G8:L12").Clear()

# Row 12 used to end with a thick bottom border because the bottom edge of
# the side tables lived there. Now that those cells have moved away, let
# the row go back to its normal (auto) height instead of the explicit
# "thick border" height.
$ws.Rows.Item(12).AutoFit()

# Restore the view: zoom level and the selected cell as left by the editor.
$ws.Application.ActiveWindow.Zoom = 160
$ws.Range("F14").Select()
